# [Documentation] Fix table names
# Rename the "PropertyValue" table reference to "ProductProperty" for the
# Spicy / Vegetarian / IngredientName / PizzaSauce rows of the
# "PizzaIngredienten bestand:" table and the Spicy / Vegetarian rows of the
# "Misc_products bestand:" table (column D), and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "PizzaIngredienten bestand:" block whose Table column (D)
# currently reads "PropertyValue" and must become "ProductProperty".
$rows = @(15, 16, 17, 18, 26, 27)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "ProductProperty"
}

# Update the remembered selection on the sheet view.
$ws.Range("E31").Select()
